$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete rows 21-27 (table now ends at row 20)
$ws.Range("A21:L27").EntireRow.Delete()

# Row 2: 에스오에스랩
$ws.Range("A2").Value = "BNK"
$ws.Range("B2").Value = "'2024-06-14"
$ws.Range("B2").Style = "Normal"
$ws.Range("C2").Value = "에스오에스랩"
$ws.Range("D2").Value = "한국"
$ws.Range("E2").Value = "한국, BNK"
$ws.Range("F2").Value = "'2024-06-19"
$ws.Range("F2").Style = "Normal"
$ws.Range("G2").Value = "'2024-06-25"
$ws.Range("G2").Style = "Normal"
$ws.Range("H2").Value = 1150
$ws.Range("I2").Value = 2000000
$ws.Range("J2").Value = 11500
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 5

# Row 3: 디비금융스팩12호
$ws.Range("A3").Value = "DB"
$ws.Range("B3").Value = "'2024-06-05"
$ws.Range("B3").Style = "Normal"
$ws.Range("C3").Value = "디비금융스팩12호"
$ws.Range("D3").Value = "DB"
$ws.Range("E3").Value = "DB"
$ws.Range("F3").Value = "'2024-06-11"
$ws.Range("F3").Style = "Normal"
$ws.Range("G3").Value = "'2024-06-18"
$ws.Range("G3").Style = "Normal"
$ws.Range("H3").Value = 10000
$ws.Range("I3").Value = 5000000
$ws.Range("J3").Value = 2000
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 100

# Row 4: 한중엔시에스
$ws.Range("A4").Value = "IBK"
$ws.Range("B4").Value = "'2024-06-10"
$ws.Range("B4").Style = "Normal"
$ws.Range("C4").Value = "한중엔시에스"
$ws.Range("D4").Value = "IBK"
$ws.Range("E4").Value = "IBK"
$ws.Range("F4").Value = "'2024-06-13"
$ws.Range("F4").Style = "Normal"
$ws.Range("G4").Value = "'2024-06-24"
$ws.Range("G4").Style = "Normal"
$ws.Range("H4").Value = 48000
$ws.Range("I4").Value = 1600000
$ws.Range("J4").Value = 30000
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 100

# Row 5: KB제28호스팩
$ws.Range("A5").Value = "KB"
$ws.Range("B5").Value = "'2024-05-07"
$ws.Range("B5").Style = "Normal"
$ws.Range("C5").Value = "KB제28호스팩"
$ws.Range("D5").Value = "KB"
$ws.Range("E5").Value = "KB"
$ws.Range("F5").Value = "'2024-05-10"
$ws.Range("F5").Style = "Normal"
$ws.Range("G5").Value = "'2024-05-17"
$ws.Range("G5").Style = "Normal"
$ws.Range("H5").Value = 10000
$ws.Range("I5").Value = 5000000
$ws.Range("J5").Value = 2000
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 100

# Row 6: KB제29호스팩
$ws.Range("A6").Value = "KB"
$ws.Range("B6").Value = "'2024-06-11"
$ws.Range("B6").Style = "Normal"
$ws.Range("C6").Value = "KB제29호스팩"
$ws.Range("D6").Value = "KB"
$ws.Range("E6").Value = "KB"
$ws.Range("F6").Value = "'2024-06-14"
$ws.Range("F6").Style = "Normal"
$ws.Range("G6").Value = "'2024-06-21"
$ws.Range("G6").Style = "Normal"
$ws.Range("H6").Value = 12000
$ws.Range("I6").Value = 6000000
$ws.Range("J6").Value = 2000
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 100

# Row 7: 에이치브이엠
$ws.Range("A7").Value = "NH"
$ws.Range("B7").Value = "'2024-06-19"
$ws.Range("B7").Style = "Normal"
$ws.Range("C7").Value = "에이치브이엠"
$ws.Range("D7").Value = "NH"
$ws.Range("E7").Value = "NH"
$ws.Range("F7").Value = "'2024-06-24"
$ws.Range("F7").Style = "Normal"
$ws.Range("G7").Value = "'2024-06-28"
$ws.Range("G7").Style = "Normal"
$ws.Range("H7").Value = 43200
$ws.Range("I7").Value = 2400000
$ws.Range("J7").Value = 18000
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 100

# Row 8: 아이씨티케이
$ws.Range("A8").Value = "NH"
$ws.Range("B8").Value = "'2024-05-07"
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").Value = "아이씨티케이"
$ws.Range("D8").Value = "NH"
$ws.Range("E8").Value = "NH"
$ws.Range("F8").Value = "'2024-05-10"
$ws.Range("F8").Style = "Normal"
$ws.Range("G8").Value = "'2024-05-17"
$ws.Range("G8").Style = "Normal"
$ws.Range("H8").Value = 39400
$ws.Range("I8").Value = 1970000
$ws.Range("J8").Value = 20000
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 100

# Row 9: 라메디텍
$ws.Range("A9").Value = "대신"
$ws.Range("B9").Value = "'2024-06-05"
$ws.Range("B9").Style = "Normal"
$ws.Range("C9").Value = "라메디텍"
$ws.Range("D9").Value = "대신"
$ws.Range("E9").Value = "대신"
$ws.Range("F9").Value = "'2024-06-11"
$ws.Range("F9").Style = "Normal"
$ws.Range("G9").Value = "'2024-06-17"
$ws.Range("G9").Style = "Normal"
$ws.Range("H9").Value = 20768
$ws.Range("I9").Value = 1298000
$ws.Range("J9").Value = 16000
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 100

# Row 10: 미래에셋비전스팩4호
$ws.Range("A10").Value = "미래"
$ws.Range("B10").Value = "'2024-05-20"
$ws.Range("B10").Style = "Normal"
$ws.Range("C10").Value = "미래에셋비전스팩4호"
$ws.Range("D10").Value = "미래"
$ws.Range("E10").Value = "미래"
$ws.Range("F10").Value = "'2024-05-23"
$ws.Range("F10").Style = "Normal"
$ws.Range("G10").Value = "'2024-05-29"
$ws.Range("G10").Style = "Normal"
$ws.Range("H10").Value = 13300
$ws.Range("I10").Value = 6650000
$ws.Range("J10").Value = 2000
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 100

# Row 11: 미래에셋비전스팩6호
$ws.Range("A11").Value = "미래"
$ws.Range("B11").Value = "'2024-06-13"
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Value = "미래에셋비전스팩6호"
$ws.Range("D11").Value = "미래"
$ws.Range("E11").Value = "미래"
$ws.Range("F11").Value = "'2024-06-18"
$ws.Range("F11").Style = "Normal"
$ws.Range("G11").Value = "'2024-06-24"
$ws.Range("G11").Style = "Normal"
$ws.Range("H11").Value = 12900
$ws.Range("I11").Value = 6450000
$ws.Range("J11").Value = 2000
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 100

# Row 12: 미래에셋비전스팩5호
$ws.Range("A12").Value = "미래"
$ws.Range("B12").Value = "'2024-06-10"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = "미래에셋비전스팩5호"
$ws.Range("D12").Value = "미래"
$ws.Range("E12").Value = "미래"
$ws.Range("F12").Value = "'2024-06-13"
$ws.Range("F12").Style = "Normal"
$ws.Range("G12").Value = "'2024-06-19"
$ws.Range("G12").Style = "Normal"
$ws.Range("H12").Value = 9500
$ws.Range("I12").Value = 4750000
$ws.Range("J12").Value = 2000
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 100

# Row 13: 그리드위즈
$ws.Range("A13").Value = "삼성"
$ws.Range("B13").Value = "'2024-06-03"
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = "그리드위즈"
$ws.Range("D13").Value = "삼성"
$ws.Range("E13").Value = "삼성"
$ws.Range("F13").Value = "'2024-06-07"
$ws.Range("F13").Style = "Normal"
$ws.Range("G13").Value = "'2024-06-14"
$ws.Range("G13").Style = "Normal"
$ws.Range("H13").Value = 56000
$ws.Range("I13").Value = 1400000
$ws.Range("J13").Value = 40000
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 100

# Row 14: 노브랜드
$ws.Range("A14").Value = "삼성"
$ws.Range("B14").Value = "'2024-05-13"
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").Value = "노브랜드"
$ws.Range("D14").Value = "삼성"
$ws.Range("E14").Value = "삼성"
$ws.Range("F14").Value = "'2024-05-17"
$ws.Range("F14").Style = "Normal"
$ws.Range("G14").Value = "'2024-05-23"
$ws.Range("G14").Style = "Normal"
$ws.Range("H14").Value = 16800
$ws.Range("I14").Value = 1200000
$ws.Range("J14").Value = 14000
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 100

# Row 15: 한국제14호스팩
$ws.Range("A15").Value = "한국"
$ws.Range("B15").Value = "'2024-06-10"
$ws.Range("B15").Style = "Normal"
$ws.Range("C15").Value = "한국제14호스팩"
$ws.Range("D15").Value = "한국"
$ws.Range("E15").Value = "한국"
$ws.Range("F15").Value = "'2024-06-13"
$ws.Range("F15").Style = "Normal"
$ws.Range("G15").Value = "'2024-06-19"
$ws.Range("G15").Style = "Normal"
$ws.Range("H15").Value = 8000
$ws.Range("I15").Value = 4000000
$ws.Range("J15").Value = 2000
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 100

# Row 16: 에스오에스랩
$ws.Range("A16").Value = "한국"
$ws.Range("B16").Value = "'2024-06-14"
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = "에스오에스랩"
$ws.Range("D16").Value = "한국"
$ws.Range("E16").Value = "한국, BNK"
$ws.Range("F16").Value = "'2024-06-19"
$ws.Range("F16").Style = "Normal"
$ws.Range("G16").Value = "'2024-06-25"
$ws.Range("G16").Style = "Normal"
$ws.Range("H16").Value = 21850
$ws.Range("I16").Value = 2000000
$ws.Range("J16").Value = 11500
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 95

# Row 17: 씨어스테크놀로지
$ws.Range("A17").Value = "한국"
$ws.Range("B17").Value = "'2024-06-10"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = "씨어스테크놀로지"
$ws.Range("D17").Value = "한국"
$ws.Range("E17").Value = "한국"
$ws.Range("F17").Value = "'2024-06-13"
$ws.Range("F17").Style = "Normal"
$ws.Range("G17").Value = "'2024-06-19"
$ws.Range("G17").Style = "Normal"
$ws.Range("H17").Value = 22100
$ws.Range("I17").Value = 1300000
$ws.Range("J17").Value = 17000
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 100

# Row 18: 한국제15호스팩
$ws.Range("A18").Value = "한국"
$ws.Range("B18").Value = "'2024-06-17"
$ws.Range("B18").Style = "Normal"
$ws.Range("C18").Value = "한국제15호스팩"
$ws.Range("D18").Value = "한국"
$ws.Range("E18").Value = "한국"
$ws.Range("F18").Value = "'2024-06-20"
$ws.Range("F18").Style = "Normal"
$ws.Range("G18").Value = "'2024-06-26"
$ws.Range("G18").Style = "Normal"
$ws.Range("H18").Value = 12500
$ws.Range("I18").Value = 6250000
$ws.Range("J18").Value = 2000
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 100

# Row 19: 하이젠알앤엠
$ws.Range("A19").Value = "한국"
$ws.Range("B19").Value = "'2024-06-18"
$ws.Range("B19").Style = "Normal"
$ws.Range("C19").Value = "하이젠알앤엠"
$ws.Range("D19").Value = "한국"
$ws.Range("E19").Value = "한국"
$ws.Range("F19").Value = "'2024-06-21"
$ws.Range("F19").Style = "Normal"
$ws.Range("G19").Value = "'2024-06-27"
$ws.Range("G19").Style = "Normal"
$ws.Range("H19").Value = 23800
$ws.Range("I19").Value = 3400000
$ws.Range("J19").Value = 7000
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 100

# Row 20: 에이치엠씨제7호스팩
$ws.Range("A20").Value = "현대차"
$ws.Range("B20").Value = "'2024-06-11"
$ws.Range("B20").Style = "Normal"
$ws.Range("C20").Value = "에이치엠씨제7호스팩"
$ws.Range("D20").Value = "현대차"
$ws.Range("E20").Value = "현대차"
$ws.Range("F20").Value = "'2024-06-14"
$ws.Range("F20").Style = "Normal"
$ws.Range("G20").Value = "'2024-06-24"
$ws.Range("G20").Style = "Normal"
$ws.Range("H20").Value = 14000
$ws.Range("I20").Value = 7000000
$ws.Range("J20").Value = 2000
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 100
